# Generate Report for Handback
# -----------------------------------------------------------------------
# This mirrors the localization CI job that refreshes localization-status.xlsx
# after a handback round-trips: the zh-cn / de-de rows grow a "Latest Target
# File" + "Latest Handback File" pointer, de-de also records the handback
# timestamp, the Overview's "Status" text flips to the handed-back message,
# and a few report columns are widened so the longer file names are legible.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatusText = "Handed back: in sync with en-US"

$mdFile      = "4a76797a-62b4-4363-a72e-772f7aa5d7e7.md"
$mdFileUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80b713d1ed3db5f77c83ecb9a2c60fa81273bf14/e2e/4a76797a-62b4-4363-a72e-772f7aa5d7e7.md"
$mdFile2     = "ffff80eb7111-9def-4022-94b6-787d58781a3d.md"
$mdFile2Url  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80b713d1ed3db5f77c83ecb9a2c60fa81273bf14/e2e/ffff80eb7111-9def-4022-94b6-787d58781a3d.md"

$zhCnHandback = "4a76797a-62b4-4363-a72e-772f7aa5d7e7.4b7ae30990be9f52f1a823ef7d37715c5d818e38.zh-cn.xlf"
$deDeHandback = "4a76797a-62b4-4363-a72e-772f7aa5d7e7.4b7ae30990be9f52f1a823ef7d37715c5d818e38.de-de.xlf"

$zhCnHandbackTime = "2016-08-24 00:59:56"
$deDeHandbackTime = "2016-08-24 01:00:17"

# widened columns (report now shows full handoff/handback file names)
$wideWidth = 29.166666666666668
$fullWidth = 39.166666666666664

# --- Overview sheet -----------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = $wideWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideWidth

$wsOverview.Range("E2").Value = $newStatusText
$wsOverview.Range("F2").Value = $newStatusText
$wsOverview.Range("E3").Value = $newStatusText
$wsOverview.Range("F3").Value = $newStatusText

# --- zh-cn sheet ----------------------------------------------------------
$wsZhCn.Columns.Item(3).ColumnWidth = $wideWidth
$wsZhCn.Columns.Item(9).ColumnWidth = $fullWidth
$wsZhCn.Columns.Item(10).ColumnWidth = $fullWidth

$wsZhCn.Range("C2").Value = $newStatusText
$wsZhCn.Range("C3").Value = $newStatusText

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdFileUrl, "", "", $mdFile) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdFileUrl, "", "", $mdFile) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $mdFile2Url, "", "", $mdFile2) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $mdFileUrl, "", "", $mdFile) | Out-Null

$wsZhCn.Range("J2").Value = $zhCnHandback
$wsZhCn.Range("J3").Value = $zhCnHandback

# --- de-de sheet ----------------------------------------------------------
$wsDeDe.Columns.Item(3).ColumnWidth = $wideWidth
$wsDeDe.Columns.Item(9).ColumnWidth = $fullWidth
$wsDeDe.Columns.Item(10).ColumnWidth = $fullWidth

$wsDeDe.Range("C2").Value = $newStatusText
$wsDeDe.Range("C3").Value = $newStatusText

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdFileUrl, "", "", $mdFile) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdFileUrl, "", "", $mdFile) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $mdFile2Url, "", "", $mdFile2) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $mdFileUrl, "", "", $mdFile) | Out-Null

$wsDeDe.Range("J2").Value = $deDeHandback
$wsDeDe.Range("J3").Value = $deDeHandback

$wsDeDe.Range("K2").Value = $deDeHandbackTime
$wsDeDe.Range("K3").Value = $deDeHandbackTime

# zh-cn's handback datetime column keeps referencing the shared "handback
# pending" string slot, whose text now reflects the zh-cn handback time.
$wsZhCn.Range("K2").Value = $zhCnHandbackTime
$wsZhCn.Range("K3").Value = $zhCnHandbackTime

Write-Host "Report regenerated for handback."
